$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.640.41"
$ws.Range("D3").Value = "2.446.61"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "2.445.73"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("B11").Value = "BabyDogeCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D11").Value = "0.0₅0199"
$ws.Range("E11").Value = "  +597.63%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.160"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("B14").Value = "Cardano"
$ws.Range("C14").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.353"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.92%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.95%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.890.42"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.532.39"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.443.05"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "326.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.73%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "646.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.80%  "
$ws.Range("E28").Value = "  +14.27%  "
$ws.Range("E29").Value = "  +6.09%  "
$ws.Range("D30").Value = "0.0₃0977"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  +6.65%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.61%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.30%  "
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +28.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "144.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.81%  "
